$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the style of the existing header cells (copy format only from H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Column I and J values per row
$valuesI = @{2=8; 3=8; 4=7; 5=8; 6=6; 7=7; 8=5}
$valuesJ = @{2=8; 3=8; 4=7; 5=8; 6=6; 7=7; 8=5}

foreach ($row in 2..8) {
    $ws.Cells.Item($row, 9).Value = $valuesI[$row]
    $ws.Cells.Item($row, 10).Value = $valuesJ[$row]
}
